# Add a new worksheet "Hoja2" placed after "Hoja1"
$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws  = $wb.Worksheets.Add($null, $ws1)
$ws.Name = "Hoja2"

# --- Cell content -----------------------------------------------------
# Written in the same order the original author typed it in: the table
# body (rows 2-16, left-to-right) first, the "work" header cell last.
$ws.Range("A2").Value  = "Campo"
$ws.Range("B2").Value  = "Columna Excel"
$ws.Range("A3").Value  = "IdFasecolda"
$ws.Range("B3").Value  = "N/A"
$ws.Range("A4").Value  = "IdControl"
$ws.Range("B4").Value  = "N/A"
$ws.Range("A5").Value  = "DatosAuto"
$ws.Range("B5").Value  = "N/A"
$ws.Range("A6").Value  = "FechaInsert"
$ws.Range("B6").Value  = "Sysdate"
$ws.Range("A7").Value  = "EstadoProceso"
$ws.Range("B7").Value  = "Por Defecto N"
$ws.Range("A8").Value  = "DescVehiculo"
$ws.Range("B8").Value  = "Referencia2 +' '+ Referencia3"
$ws.Range("A9").Value  = "DesVehiculoTipo"
$ws.Range("B9").Value  = "Clase"
$ws.Range("A10").Value = "CodVehiculoTipo"
$ws.Range("B10").Value = "N/A"
$ws.Range("A11").Value = "DesVehiculoMarca"
$ws.Range("B11").Value = "Marca"
$ws.Range("A12").Value = "CodVehiculoMarca"
$ws.Range("B12").Value = "N/A"
$ws.Range("A13").Value = "CodVehiculoExterno"
$ws.Range("B13").Value = "Codigo"
$ws.Range("A14").Value = "CodHomologado"
$ws.Range("B14").Value = "Homologocodigo"
$ws.Range("A15").Value = "DesGrupoModelo"
$ws.Range("B15").Value = "Referencia1"
$ws.Range("A16").Value = "CodGruproModelo"
$ws.Range("B16").Value = "N/A"
$ws.Range("B1").Value  = "work"

# --- Column widths ------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 19.17
$ws.Columns.Item(2).ColumnWidth = 28.59

# --- Styling --------------------------------------------------------
# Helper: build each distinct look once on a scratch cell, copy it, then
# paste-special (formats only) onto every cell that shares it so the
# style table doesn't balloon with duplicate entries.
$xlPasteFormats = -4122
$xlContinuous   = 1
$xlCenter       = -4108
$white          = 16777215
$black          = 0
$darkGray       = 0x222222
$edgeLeft   = 7
$edgeTop    = 8
$edgeBottom = 9
$edgeRight  = 10
$medium     = -4138

# Style for A2 (top-left header cell): bold 12pt black Aptos Narrow,
# white fill, box border on all sides, vertically centered.
$s1 = $ws.Range("Z1")
$s1.Font.Bold = $true
$s1.Font.Size = 12
$s1.Font.Color = $black
$s1.Font.Name = "Aptos Narrow"
$s1.Interior.Color = $white
$s1.VerticalAlignment = $xlCenter
$s1.Borders.Item($edgeLeft).LineStyle = $xlContinuous
$s1.Borders.Item($edgeLeft).Weight = $medium
$s1.Borders.Item($edgeRight).LineStyle = $xlContinuous
$s1.Borders.Item($edgeRight).Weight = $medium
$s1.Borders.Item($edgeTop).LineStyle = $xlContinuous
$s1.Borders.Item($edgeTop).Weight = $medium
$s1.Borders.Item($edgeBottom).LineStyle = $xlContinuous
$s1.Borders.Item($edgeBottom).Weight = $medium
$s1.Copy()
$ws.Range("A2").PasteSpecial($xlPasteFormats)
$s1.Clear()

# Style for B2 (top-right header cell): same font/fill, but no left edge.
$s2 = $ws.Range("Z1")
$s2.Font.Bold = $true
$s2.Font.Size = 12
$s2.Font.Color = $black
$s2.Font.Name = "Aptos Narrow"
$s2.Interior.Color = $white
$s2.VerticalAlignment = $xlCenter
$s2.Borders.Item($edgeRight).LineStyle = $xlContinuous
$s2.Borders.Item($edgeRight).Weight = $medium
$s2.Borders.Item($edgeTop).LineStyle = $xlContinuous
$s2.Borders.Item($edgeTop).Weight = $medium
$s2.Borders.Item($edgeBottom).LineStyle = $xlContinuous
$s2.Borders.Item($edgeBottom).Weight = $medium
$s2.Copy()
$ws.Range("B2").PasteSpecial($xlPasteFormats)
$s2.Clear()

# Style for column A, rows 3-16: regular 12pt black Aptos Narrow, white
# fill, left+right+bottom border.
$s3 = $ws.Range("Z1")
$s3.Font.Size = 12
$s3.Font.Color = $black
$s3.Font.Name = "Aptos Narrow"
$s3.Interior.Color = $white
$s3.VerticalAlignment = $xlCenter
$s3.Borders.Item($edgeLeft).LineStyle = $xlContinuous
$s3.Borders.Item($edgeLeft).Weight = $medium
$s3.Borders.Item($edgeRight).LineStyle = $xlContinuous
$s3.Borders.Item($edgeRight).Weight = $medium
$s3.Borders.Item($edgeBottom).LineStyle = $xlContinuous
$s3.Borders.Item($edgeBottom).Weight = $medium
$s3.Copy()
$ws.Range("A3:A16").PasteSpecial($xlPasteFormats)
$s3.Clear()

# Style for column B, "N/A" / plain rows: same font as s3, but only
# right+bottom border.
$s4 = $ws.Range("Z1")
$s4.Font.Size = 12
$s4.Font.Color = $black
$s4.Font.Name = "Aptos Narrow"
$s4.Interior.Color = $white
$s4.VerticalAlignment = $xlCenter
$s4.Borders.Item($edgeRight).LineStyle = $xlContinuous
$s4.Borders.Item($edgeRight).Weight = $medium
$s4.Borders.Item($edgeBottom).LineStyle = $xlContinuous
$s4.Borders.Item($edgeBottom).Weight = $medium
$s4.Copy()
$ws.Range("B3:B7").PasteSpecial($xlPasteFormats)
$ws.Range("B10").PasteSpecial($xlPasteFormats)
$s4.Clear()

# Style for column B, "web paste" rows: bold 12pt dark-gray Calibri,
# right+bottom border.
$s5 = $ws.Range("Z1")
$s5.Font.Bold = $true
$s5.Font.Size = 12
$s5.Font.Color = $darkGray
$s5.Font.Name = "Calibri"
$s5.Interior.Color = $white
$s5.VerticalAlignment = $xlCenter
$s5.Borders.Item($edgeRight).LineStyle = $xlContinuous
$s5.Borders.Item($edgeRight).Weight = $medium
$s5.Borders.Item($edgeBottom).LineStyle = $xlContinuous
$s5.Borders.Item($edgeBottom).Weight = $medium
$s5.Copy()
$ws.Range("B8:B9").PasteSpecial($xlPasteFormats)
$ws.Range("B11:B14").PasteSpecial($xlPasteFormats)
$s5.Clear()

# Style for B15:B16: same as s1/s2 font (bold black Aptos Narrow), right+bottom border.
$s6 = $ws.Range("Z1")
$s6.Font.Bold = $true
$s6.Font.Size = 12
$s6.Font.Color = $black
$s6.Font.Name = "Aptos Narrow"
$s6.Interior.Color = $white
$s6.VerticalAlignment = $xlCenter
$s6.Borders.Item($edgeRight).LineStyle = $xlContinuous
$s6.Borders.Item($edgeRight).Weight = $medium
$s6.Borders.Item($edgeBottom).LineStyle = $xlContinuous
$s6.Borders.Item($edgeBottom).Weight = $medium
$s6.Copy()
$ws.Range("B15:B16").PasteSpecial($xlPasteFormats)
$s6.Clear()

# --- Row heights ------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 15.75
$ws.Rows.Item(2).RowHeight = 16.5
$ws.Rows.Item(3).RowHeight = 16.5
$ws.Rows.Item(4).RowHeight = 16.5
$ws.Rows.Item(5).RowHeight = 16.5
$ws.Rows.Item(6).RowHeight = 16.5
$ws.Rows.Item(7).RowHeight = 16.5
$ws.Rows.Item(8).RowHeight = 16.5
$ws.Rows.Item(9).RowHeight = 16.5
$ws.Rows.Item(10).RowHeight = 16.5
$ws.Rows.Item(11).RowHeight = 16.5
$ws.Rows.Item(12).RowHeight = 16.5
$ws.Rows.Item(13).RowHeight = 16.5
$ws.Rows.Item(14).RowHeight = 16.5
$ws.Rows.Item(15).RowHeight = 16.5
$ws.Rows.Item(16).RowHeight = 16.5

# --- Selections / active sheet ----------------------------------------
# Hoja1: selection moves to C6, it's no longer the active/visible tab.
$ws1.Range("C6").Select()

# Hoja2: selection on E16, and it becomes the active (visible) tab.
$ws.Range("E16").Select()
$ws.Activate()
